# This edit inserts one new data row (for a new "Primera" quality Mango
# price observation dated 44572) above the existing row 203 on Sheet1.
# Inserting the row pushes the existing rows 203-314 down to 204-315,
# and Excel carries the row formatting (e.g. the date style on column D)
# down with them. We then populate the newly inserted row 203 with its
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 203; everything below (203-314) shifts to 204-315.
$ws.Rows(203).Insert()

# Fill in the values for the newly inserted row 203.
$ws.Range("A203").Value = 3
$ws.Range("B203").Value = "Femacal de La Calera"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = 44572
$ws.Range("E203").Value = 5
$ws.Range("F203").Value = "Fruta"
$ws.Range("G203").Value = 100108
$ws.Range("H203").Value = "Tropicales y subtropicales"
$ws.Range("I203").Value = 100108002
$ws.Range("J203").Value = "Mango"
$ws.Range("K203").Value = "Sin especificar"
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 456
$ws.Range("N203").Value = 7000
$ws.Range("O203").Value = 7000
$ws.Range("P203").Value = 7000
$ws.Range("Q203").Value = "`$/bandeja 4 kilos"
$ws.Range("R203").Value = "Perú"
$ws.Range("S203").Value = 1750
$ws.Range("T203").Value = 4
